# Update with restock suggestion
# - Remove the "Sales Volume Rank" column (old Q) from "Forecast Comparison";
#   "Lifecycle Stage" shifts left into column Q.
# - Populate Week_Start_Date (col B) for every data row.
# - Refresh Inventory Coverage (L), Seasonality Index (P) and the
#   (now-shifted) Lifecycle Stage (Q, values become "Decline") for every row.
# - A few MyForecast (D) values changed slightly.
# - Refresh a handful of rollup figures on the Summary sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Drop the old "Sales Volume Rank" column outright; everything to its right
# (just "Lifecycle Stage") shifts left to take its place.
$ws.Range("Q:Q").Delete()

# --- Week_Start_Date (column B) ----------------------------------------
# Force these to be stored as plain text (not auto-converted to date
# serials) by marking the range as Text before writing the values.
$ws.Range("B2:B17").NumberFormat = "@"

$weekStartDates = @{
    2  = "2025-02-02"
    3  = "2025-02-09"
    4  = "2025-02-16"
    5  = "2025-02-23"
    6  = "2025-03-02"
    7  = "2025-03-09"
    8  = "2025-03-16"
    9  = "2025-03-23"
    10 = "2025-03-30"
    11 = "2025-04-06"
    12 = "2025-04-13"
    13 = "2025-04-20"
    14 = "2025-04-27"
    15 = "2025-05-04"
    16 = "2025-05-11"
    17 = "2025-05-18"
}
foreach ($row in $weekStartDates.Keys) {
    $ws.Range("B$row").Value = $weekStartDates[$row]
}

# --- MyForecast (column D) small corrections ----------------------------
$ws.Range("D9").Value = 118
$ws.Range("D14").Value = 103
$ws.Range("D15").Value = 97

# --- Inventory Coverage (column L) --------------------------------------
$inventoryCoverage = @{
    2  = 0.34
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
}
foreach ($row in $inventoryCoverage.Keys) {
    $ws.Range("L$row").Value = $inventoryCoverage[$row]
}

# --- Seasonality Index (column P) ---------------------------------------
$seasonalityIndex = @{
    2  = 0.87
    3  = 1.12
    4  = 1.16
    5  = 0.95
    6  = 1.1
    7  = 1.1
    8  = 1.07
    9  = 1.11
    10 = 0.8
    11 = 0.94
    12 = 0.88
    13 = 1.07
    14 = 0.94
    15 = 1.03
    16 = 0.88
    17 = 1.2
}
foreach ($row in $seasonalityIndex.Keys) {
    $ws.Range("P$row").Value = $seasonalityIndex[$row]
}

# --- Lifecycle Stage (now column Q after the delete) ---------------------
for ($row = 2; $row -le 17; $row++) {
    $ws.Range("Q$row").Value = "Decline"
}

# --- Summary sheet rollups -------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("B9:B10").NumberFormat = "@"
$ws2.Range("B9").Value = "1678"
$ws2.Range("B10").Value = "856"
$ws2.Range("B13").Value = "N/A"
$ws2.Range("B15").Value = "N/A"

Write-Output "Applied restock-suggestion update"
